# Comercializadora del Agro de Limarí - Poroto verde : weekly price update
#
# The new data snapshot inserts 3 new daily-price records right after the
# report header/summary rows (rows 2-6) and before the existing detail
# rows, which pushes every existing detail row (old rows 7-105) down by
# three positions (they land on new rows 10-108). Because the old rows
# simply slide down as a block, a single "insert 3 rows" operation
# reproduces that whole shift (including the 3 extra rows that appear at
# the bottom of the sheet, 106-108, which end up holding what used to be
# rows 103-105) -- so all that's left by hand is filling in the 3 brand
# new rows with their own data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push existing rows 7-105 down to 10-108, leaving 3 blank rows (7,8,9).
$ws.Range("A7:A9").EntireRow.Insert()

# Row 7: new "Magnum / Primera" record dated 2021-11-03
$ws.Cells.Item(7, 1).Value = 2
$ws.Cells.Item(7, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(7, 3).Value = "Coquimbo"
$ws.Cells.Item(7, 4).Value = 44503
$ws.Cells.Item(7, 5).Value = 4
$ws.Cells.Item(7, 6).Value = 100112031
$ws.Cells.Item(7, 7).Value = "Poroto verde"
$ws.Cells.Item(7, 8).Value = "Magnum"
$ws.Cells.Item(7, 9).Value = "Primera"
$ws.Cells.Item(7, 10).Value = 500
$ws.Cells.Item(7, 11).Value = 33000
$ws.Cells.Item(7, 12).Value = 33000
$ws.Cells.Item(7, 13).Value = 33000
$ws.Cells.Item(7, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(7, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(7, 16).Value = 1320
$ws.Cells.Item(7, 17).Value = 25
$ws.Cells.Item(7, 18).Value = "Hortaliza"

# Row 8: new "Magnum / Segunda" record dated 2021-11-03
$ws.Cells.Item(8, 1).Value = 2
$ws.Cells.Item(8, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(8, 3).Value = "Coquimbo"
$ws.Cells.Item(8, 4).Value = 44503
$ws.Cells.Item(8, 5).Value = 4
$ws.Cells.Item(8, 6).Value = 100112031
$ws.Cells.Item(8, 7).Value = "Poroto verde"
$ws.Cells.Item(8, 8).Value = "Magnum"
$ws.Cells.Item(8, 9).Value = "Segunda"
$ws.Cells.Item(8, 10).Value = 500
$ws.Cells.Item(8, 11).Value = 31000
$ws.Cells.Item(8, 12).Value = 31000
$ws.Cells.Item(8, 13).Value = 31000
$ws.Cells.Item(8, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(8, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(8, 16).Value = 1240
$ws.Cells.Item(8, 17).Value = 25
$ws.Cells.Item(8, 18).Value = "Hortaliza"

# Row 9: new "Sin especificar / Primera" record dated 2021-11-03
$ws.Cells.Item(9, 1).Value = 2
$ws.Cells.Item(9, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(9, 3).Value = "Coquimbo"
$ws.Cells.Item(9, 4).Value = 44503
$ws.Cells.Item(9, 5).Value = 4
$ws.Cells.Item(9, 6).Value = 100112031
$ws.Cells.Item(9, 7).Value = "Poroto verde"
$ws.Cells.Item(9, 8).Value = "Sin especificar"
$ws.Cells.Item(9, 9).Value = "Primera"
$ws.Cells.Item(9, 10).Value = 400
$ws.Cells.Item(9, 11).Value = 40000
$ws.Cells.Item(9, 12).Value = 42000
$ws.Cells.Item(9, 13).Value = 41000
$ws.Cells.Item(9, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(9, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(9, 16).Value = 1640
$ws.Cells.Item(9, 17).Value = 25
$ws.Cells.Item(9, 18).Value = "Hortaliza"
